# Update cryptocurrency Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "51.917.24"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "2.788.02"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("E4").Value = "  -0.01%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "357.72"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +1.62%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "109.15"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -3.26%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.565"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -0.96%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "40.00"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  -3.40%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0854"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("E13").Value = "  -1.98%  "

$ws.Range("E14").Value = "  -1.98%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "3.225.52"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -0.90%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "2.769.65"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  -1.96%  "

$ws.Range("E17").Value = "  +6.66%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "51.853.24"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +0.52%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "7.41"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  -0.88%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "3.13"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -1.95%  "

$ws.Range("E21").Value = "  -2.39%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0980"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -1.09%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "274.27"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +1.38%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "70.26"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("E28").Value = "  -1.25%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "0.145"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +4.12%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "2.21"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -1.56%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0468"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  +4.51%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "51.55"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +1.96%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "34.35"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +1.14%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "5.72"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -1.55%  "

$ws.Range("E35").Value = "  +2.54%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "5.26"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +1.05%  "

$ws.Range("E37").Value = "  +0.04%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "3.22"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  +0.63%  "

$ws.Range("E39").Value = "  -2.80%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "17.97"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "2.55"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("E42").Value = "  -1.59%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "2.25"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -1.52%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "121.78"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  -3.44%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "21.97"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -7.78%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "2.074.49"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -0.16%  "

$ws.Range("E47").Value = "  -2.44%  "

$ws.Range("E48").Value = "  -4.25%  "

$ws.Range("E49").Value = "  +1.38%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "0.930"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "8.92"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
